$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues = -4163
$xlPasteValues = -4163
$stageRow = 1000

function Swap-Rows($rowA, $rowB, $stage) {
    $rangeA = "A" + $rowA + ":AY" + $rowA
    $rangeB = "A" + $rowB + ":AY" + $rowB
    $rangeStage = "A" + $stage + ":AY" + $stage

    $ws.Range($rangeA).Copy()
    $ws.Range($rangeStage).PasteSpecial($xlPasteValues)

    $ws.Range($rangeA).ClearContents()
    $ws.Range($rangeB).Copy()
    $ws.Range($rangeA).PasteSpecial($xlPasteValues)

    $ws.Range($rangeB).ClearContents()
    $ws.Range($rangeStage).Copy()
    $ws.Range($rangeB).PasteSpecial($xlPasteValues)

    $ws.Range($rangeStage).Clear()
}

Swap-Rows 6 7 $stageRow
Swap-Rows 8 9 $stageRow
Swap-Rows 15 16 $stageRow

$excel.CutCopyMode = 0
